$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.4767983333333334
$ws.Range("H2").Value = 1.430395
$ws.Range("I2").Value = 0.003723890400117776
$ws.Range("J2").Value = 0.003723890400117776
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.404869
$ws.Range("N2").Value = 4.214607
$ws.Range("O2").Value = 0.6692718564235921
$ws.Range("P2").Value = 0.6692718564235923
$ws.Range("Q2").Value = 0.6698391977516667
$ws.Range("R2").Value = 6.028552779765
$ws.Range("S2").Value = 0.002492295041204817
$ws.Range("T2").Value = 0.002492295041204818

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.4767983333333334
$ws.Range("H3").Value = 1.430395
$ws.Range("I3").Value = 0.003723890400117776
$ws.Range("J3").Value = 0.003723890400117776
$ws.Range("O3").Value = 0.3150411080808892
$ws.Range("P3").Value = 0.3150411080808893
$ws.Range("Q3").Value = 0.3153081682283334
$ws.Range("R3").Value = 2.837773514055
$ws.Range("S3").Value = 0.00117317855802489
$ws.Range("T3").Value = 0.00117317855802489

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.4767983333333334
$ws.Range("H4").Value = 1.430395
$ws.Range("I4").Value = 0.003723890400117776
$ws.Range("J4").Value = 0.003723890400117776
$ws.Range("M4").Value = 0.03292866666666667
$ws.Range("N4").Value = 0.098786
$ws.Range("O4").Value = 0.01568703549551856
$ws.Range("P4").Value = 0.01568703549551856
$ws.Range("Q4").Value = 0.01570033338555556
$ws.Range("R4").Value = 0.14130300047
$ws.Range("S4").Value = 0.00005841680088806837
$ws.Range("T4").Value = 0.00005841680088806837

$ws.Range("I5").Value = 0.8490200321922391
$ws.Range("J5").Value = 0.8490200321922391
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.404869
$ws.Range("N5").Value = 4.214607
$ws.Range("O5").Value = 0.6692718564235921
$ws.Range("P5").Value = 0.6692718564235923
$ws.Range("Q5").Value = 152.7184842015643
$ws.Range("R5").Value = 1374.466357814079
$ws.Range("S5").Value = 0.5682252130861178
$ws.Range("T5").Value = 0.568225213086118

$ws.Range("I6").Value = 0.8490200321922391
$ws.Range("J6").Value = 0.8490200321922391
$ws.Range("O6").Value = 0.3150411080808892
$ws.Range("P6").Value = 0.3150411080808893
$ws.Range("S6").Value = 0.2674762117247152
$ws.Range("T6").Value = 0.2674762117247153

$ws.Range("I7").Value = 0.8490200321922391
$ws.Range("J7").Value = 0.8490200321922391
$ws.Range("M7").Value = 0.03292866666666667
$ws.Range("N7").Value = 0.098786
$ws.Range("O7").Value = 0.01568703549551856
$ws.Range("P7").Value = 0.01568703549551856
$ws.Range("Q7").Value = 3.579562265315778
$ws.Range("R7").Value = 32.216060387842
$ws.Range("S7").Value = 0.01331860738140596
$ws.Range("T7").Value = 0.01331860738140596

$ws.Range("G8").Value = 18.85432833333333
$ws.Range("H8").Value = 56.562985
$ws.Range("I8").Value = 0.1472560774076432
$ws.Range("J8").Value = 0.1472560774076432
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.404869
$ws.Range("N8").Value = 4.214607
$ws.Range("O8").Value = 0.6692718564235921
$ws.Range("P8").Value = 0.6692718564235923
$ws.Range("Q8").Value = 26.48786139132166
$ws.Range("R8").Value = 238.390752521895
$ws.Range("S8").Value = 0.09855434829626951
$ws.Range("T8").Value = 0.09855434829626956

$ws.Range("G9").Value = 18.85432833333333
$ws.Range("H9").Value = 56.562985
$ws.Range("I9").Value = 0.1472560774076432
$ws.Range("J9").Value = 0.1472560774076432
$ws.Range("O9").Value = 0.3150411080808892
$ws.Range("P9").Value = 0.3150411080808893
$ws.Range("Q9").Value = 12.46842388981833
$ws.Range("R9").Value = 112.215815008365
$ws.Range("S9").Value = 0.04639171779814909
$ws.Range("T9").Value = 0.04639171779814911

$ws.Range("G10").Value = 18.85432833333333
$ws.Range("H10").Value = 56.562985
$ws.Range("I10").Value = 0.1472560774076432
$ws.Range("J10").Value = 0.1472560774076432
$ws.Range("M10").Value = 0.03292866666666667
$ws.Range("N10").Value = 0.098786
$ws.Range("O10").Value = 0.01568703549551856
$ws.Range("P10").Value = 0.01568703549551856
$ws.Range("Q10").Value = 0.6208478929122222
$ws.Range("R10").Value = 5.587631036209999
$ws.Range("S10").Value = 0.002310011313224527
$ws.Range("T10").Value = 0.002310011313224527
